$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6100
$ws.Range("I40").Value = 12250
$ws.Range("J40").Value = 5153.846
$ws.Range("K40").Value = 12250
$ws.Range("L40").Value = 5153.846
$ws.Range("M40").Value = -12075
$ws.Range("N40").Value = -5503.846
$ws.Range("H100").Value = 2375.889
$ws.Range("I100").Value = 1769.1428
$ws.Range("K100").Value = 1769.1428
$ws.Range("M100").Value = -1228.1428
$ws.Range("H113").Value = 5188.136
$ws.Range("I113").Value = 5009.0713
$ws.Range("J113").Value = 5501.5
$ws.Range("K113").Value = 5009.0713
$ws.Range("L113").Value = 5501.5
$ws.Range("M113").Value = -1755.0713
$ws.Range("N113").Value = -12009.5
$ws.Range("H115").Value = 1999.5
$ws.Range("I115").Value = 1999.5
$ws.Range("K115").Value = 5998.5
$ws.Range("M115").Value = -4431.5
$ws.Range("H116").Value = 10749.75
$ws.Range("I116").Value = 10166.333
$ws.Range("K116").Value = 10166.333
$ws.Range("M116").Value = -6724.333000000001
$ws.Range("H132").Value = 25002482
$ws.Range("I132").Value = 25643468
$ws.Range("K132").Value = 76930404
$ws.Range("M132").Value = -76927874
$ws.Range("I135").Value = 568.5
$ws.Range("K135").Value = 5116.5
$ws.Range("M135").Value = -2581.5
$ws.Range("H137").Value = 10884.786
$ws.Range("I137").Value = 2638.6667
$ws.Range("K137").Value = 7916.000100000001
$ws.Range("M137").Value = -5366.000100000001
$ws.Range("H141").Value = 9502.77
$ws.Range("I141").Value = 13482
$ws.Range("J141").Value = 3136
$ws.Range("K141").Value = 40446
$ws.Range("L141").Value = 9408
$ws.Range("M141").Value = -35266
$ws.Range("N141").Value = -19768

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12347480
$ws.Range("I32").Value = 13890402
$ws.Range("K32").Value = 13890402
$ws.Range("M32").Value = -13890115
$ws.Range("H45").Value = 5028
$ws.Range("I45").Value = 4174.25
$ws.Range("J45").Value = 6166.3335
$ws.Range("K45").Value = 4174.25
$ws.Range("L45").Value = 6166.3335
$ws.Range("M45").Value = -3797.25
$ws.Range("N45").Value = -6920.3335
$ws.Range("H74").Value = 21741608
$ws.Range("J74").Value = 2827.4
$ws.Range("L74").Value = 2827.4
$ws.Range("N74").Value = -4575.4
$ws.Range("H77").Value = 21741608
$ws.Range("J77").Value = 2827.4
$ws.Range("L77").Value = 14137
$ws.Range("N77").Value = -22873
$ws.Range("H132").Value = 22254744
$ws.Range("I132").Value = 1573.7407
$ws.Range("K132").Value = 4721.2221
$ws.Range("M132").Value = -2191.2221

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6526.5557
$ws.Range("I86").Value = 8350.888999999999
$ws.Range("J86").Value = 2877.889
$ws.Range("K86").Value = 8350.888999999999
$ws.Range("L86").Value = 2877.889
$ws.Range("M86").Value = -7227.888999999999
$ws.Range("N86").Value = -5123.889
$ws.Range("H89").Value = 6526.5557
$ws.Range("I89").Value = 8350.888999999999
$ws.Range("J89").Value = 2877.889
$ws.Range("K89").Value = 41754.44499999999
$ws.Range("L89").Value = 14389.445
$ws.Range("M89").Value = -36138.44499999999
$ws.Range("N89").Value = -25621.445
$ws.Range("H134").Value = 1657.3928
$ws.Range("I134").Value = 1781.5435
$ws.Range("K134").Value = 5344.6305
$ws.Range("M134").Value = -2809.6305

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 83340240
$ws.Range("I31").Value = 5312.7144
$ws.Range("K31").Value = 5312.7144
$ws.Range("M31").Value = -5017.7144
$ws.Range("H34").Value = 83340240
$ws.Range("I34").Value = 5312.7144
$ws.Range("K34").Value = 5312.7144
$ws.Range("M34").Value = -5110.7144
$ws.Range("H132").Value = 2435.7031
$ws.Range("I132").Value = 2370.9827
$ws.Range("J132").Value = 3061.3333
$ws.Range("K132").Value = 7112.9481
$ws.Range("L132").Value = 9183.999899999999
$ws.Range("M132").Value = -4582.9481
$ws.Range("N132").Value = -14243.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 23264622
$ws.Range("I4").Value = 51945416
$ws.Range("K4").Value = 155836248
$ws.Range("M4").Value = -155836136
$ws.Range("H92").Value = 799.4
$ws.Range("I92").Value = 750
$ws.Range("J92").Value = 873.5
$ws.Range("K92").Value = 2250
$ws.Range("L92").Value = 2620.5
$ws.Range("M92").Value = -1002
$ws.Range("N92").Value = -5116.5
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H117").Value = 1115
$ws.Range("J117").Value = 2000
$ws.Range("L117").Value = 6000
$ws.Range("N117").Value = -12884
$ws.Range("H120").Value = 995
$ws.Range("I120").Value = 995
$ws.Range("K120").Value = 2985
$ws.Range("M120").Value = 1853
$ws.Range("H128").Value = 189660
$ws.Range("I128").Value = 189660
$ws.Range("K128").Value = 568980
$ws.Range("M128").Value = -564000
$ws.Range("H130").Value = 1878.2
$ws.Range("I130").Value = 1108.3334
$ws.Range("K130").Value = 3325.0002
$ws.Range("M130").Value = 1694.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8875
$ws.Range("I80").Value = 12333
$ws.Range("J80").Value = 6800.2
$ws.Range("K80").Value = 12333
$ws.Range("L80").Value = 6800.2
$ws.Range("M80").Value = -11335
$ws.Range("N80").Value = -8796.200000000001
$ws.Range("H83").Value = 8875
$ws.Range("I83").Value = 12333
$ws.Range("J83").Value = 6800.2
$ws.Range("K83").Value = 61665
$ws.Range("L83").Value = 34001
$ws.Range("M83").Value = -56673
$ws.Range("N83").Value = -43985
$ws.Range("H102").Value = 4481.478
$ws.Range("I102").Value = 3391.8572
$ws.Range("K102").Value = 3391.8572
$ws.Range("M102").Value = -1769.8572
$ws.Range("H107").Value = 524.8929000000001
$ws.Range("J107").Value = 406.5
$ws.Range("L107").Value = 406.5
$ws.Range("N107").Value = -4246.5
$ws.Range("H113").Value = 3711.6875
$ws.Range("I113").Value = 3278.4119
$ws.Range("J113").Value = 4202.7334
$ws.Range("K113").Value = 3278.4119
$ws.Range("L113").Value = 4202.7334
$ws.Range("M113").Value = -1108.4119
$ws.Range("N113").Value = -8542.733400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2549.8572
$ws.Range("I22").Value = 2570.8572
$ws.Range("K22").Value = 2570.8572
$ws.Range("M22").Value = -2275.8572
$ws.Range("H27").Value = 2549.8572
$ws.Range("I27").Value = 2570.8572
$ws.Range("K27").Value = 2570.8572
$ws.Range("M27").Value = -2463.8572
$ws.Range("H46").Value = 1830.8928
$ws.Range("I46").Value = 990.4545000000001
$ws.Range("J46").Value = 4912.5
$ws.Range("K46").Value = 990.4545000000001
$ws.Range("L46").Value = 4912.5
$ws.Range("M46").Value = -802.4545000000001
$ws.Range("N46").Value = -5288.5
$ws.Range("H68").Value = 3539.8
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 4566.3335
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 4566.3335
$ws.Range("M68").Value = -1251
$ws.Range("N68").Value = -6064.3335
$ws.Range("H71").Value = 3539.8
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 4566.3335
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 22831.6675
$ws.Range("M71").Value = -6256
$ws.Range("N71").Value = -30319.6675
$ws.Range("H82").Value = 2219
$ws.Range("I82").Value = 1843.7142
$ws.Range("K82").Value = 1843.7142
$ws.Range("M82").Value = -1482.7142
$ws.Range("H85").Value = 2219
$ws.Range("I85").Value = 1843.7142
$ws.Range("K85").Value = 1843.7142
$ws.Range("M85").Value = -595.7141999999999
$ws.Range("H131").Value = 79149.42999999999
$ws.Range("J131").Value = 89750
$ws.Range("L131").Value = 89750
$ws.Range("N131").Value = -99830
$ws.Range("H133").Value = 73442
$ws.Range("J133").Value = 73442
$ws.Range("L133").Value = 73442
$ws.Range("N133").Value = -78502
$ws.Range("H136").Value = 1002778.25
$ws.Range("I136").Value = 1335651.1
$ws.Range("J136").Value = 4159.6
$ws.Range("K136").Value = 4006953.3
$ws.Range("L136").Value = 12478.8
$ws.Range("M136").Value = -4004403.3
$ws.Range("N136").Value = -17578.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2009.05
$ws.Range("I107").Value = 1122
$ws.Range("J107").Value = 2600.4167
$ws.Range("K107").Value = 3366
$ws.Range("L107").Value = 7801.250100000001
$ws.Range("M107").Value = -1446
$ws.Range("N107").Value = -11641.2501
$ws.Range("H122").Value = 133339710
$ws.Range("I122").Value = 250002220
$ws.Range("J122").Value = 40009700
$ws.Range("K122").Value = 750006660
$ws.Range("L122").Value = 120029100
$ws.Range("M122").Value = -750004210
$ws.Range("N122").Value = -120034000
$ws.Range("H132").Value = 1136.9131
$ws.Range("I132").Value = 1008.64703
$ws.Range("K132").Value = 3025.94109
$ws.Range("M132").Value = -495.9410899999998
$ws.Range("H133").Value = 64500.8
$ws.Range("J133").Value = 64500.8
$ws.Range("L133").Value = 64500.8
$ws.Range("N133").Value = -74620.8
